$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.695.29"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "1.896.70"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'239.45"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.4799"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "'0.06536"
$ws.Range("D10").Value = "1.954.62"
$ws.Range("E10").Value = "  +5.24%  "
$ws.Range("D11").Value = "'0.07482"
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "'87.95"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "'0.6670"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").Value = "30.676.25"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "'13.30"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "2.199.40"
$ws.Range("E19").Value = "  +4.12%  "
$ws.Range("D20").Value = "'0.000007609"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'228.73"
$ws.Range("E21").Value = "  +4.71%  "
$ws.Range("D22").Value = "'5.302"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'6.217"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").Value = "'168.61"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("D26").Value = "'9.270"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "'18.60"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "'1.949"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("D29").Value = "'1.402"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").Value = "'0.09941"
$ws.Range("E30").Value = "  +8.78%  "
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").Value = "'4.015"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "'0.05053"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value = "'1.221"
$ws.Range("E34").Value = "  +7.16%  "
$ws.Range("D35").Value = "'0.7532"
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").Value = "'2.712"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "'0.01873"
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("D38").Value = "'2.653"
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("D39").Value = "'0.9207"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").Value = "'2.076"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").Value = "'107.01"
$ws.Range("D42").Value = "'5.846"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "'0.4292"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("D44").Value = "'1.005"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "'7.385"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").Value = "'64.43"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("E48").Value = "  -5.31%  "
$ws.Range("D49").Value = "'8.967"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").Value = "'33.88"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").Value = "  -0.60%  "
